$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# (e.g. "242.05") are not auto-converted to number cells, matching the
# original inlineStr cell type in the workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '96.472.55'
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('D3').Value = '3.587.01'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '242.05'
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('D6').Value = '654.51'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').Value = '1.57'
$ws.Range('E7').Value = '  +6.72%  '
$ws.Range('D8').Value = '0.407'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  +3.88%  '
$ws.Range('D11').Value = '3.586.70'
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').Value = '43.31'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').Value = '6.40'
$ws.Range('E14').Value = '  +1.49%  '
$ws.Range('D15').Value = '4.251.69'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').Value = '96.281.28'
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('D17').Value = '0.0000260'
$ws.Range('E17').Value = '  +1.61%  '
$ws.Range('D18').Value = '3.590.49'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '7.77'
$ws.Range('E19').Value = '  -5.16%  '
$ws.Range('D20').Value = '12.54'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').Value = '17.76'
$ws.Range('E21').Value = '  -1.85%  '
$ws.Range('D22').Value = '0.495'
$ws.Range('E22').Value = '  +2.23%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = '512.32'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').Value = '3.44'
$ws.Range('E24').Value = '  -2.50%  '
$ws.Range('D25').Value = '0.0000204'
$ws.Range('E25').Value = '  +3.98%  '
$ws.Range('D26').Value = '6.84'
$ws.Range('E26').Value = '  +3.13%  '
$ws.Range('D27').Value = '96.44'
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('D28').Value = '12.76'
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').Value = '3.779.07'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('D30').Value = '2.98'
$ws.Range('E30').Value = '  -6.55%  '
$ws.Range('D31').Value = '0.149'
$ws.Range('E31').Value = '  +7.42%  '
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').Value = '  +3.94%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('D36').Value = '31.68'
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('D37').Value = '617.55'
$ws.Range('E37').Value = '  +8.68%  '
$ws.Range('D38').Value = '8.69'
$ws.Range('E38').Value = '  +4.89%  '
$ws.Range('E39').Value = '  +1.28%  '
$ws.Range('D40').Value = '1.63'
$ws.Range('E40').Value = '  +9.35%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').Value = '0.909'
$ws.Range('E43').Value = '  -2.06%  '
$ws.Range('D44').Value = '1.82'
$ws.Range('E44').Value = '  +5.73%  '
$ws.Range('D45').Value = '5.72'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').Value = '2.29'
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').Value = '34.28'
$ws.Range('E47').Value = '  +1.93%  '
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('D49').Value = '0.0418'
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('D50').Value = '3.58'
$ws.Range('E50').Value = '  +3.94%  '
$ws.Range('E51').Value = '  +3.03%  '

# Restore the default style on column D so no stray style index is left
# behind from the temporary text NumberFormat applied above.
$ws.Range("D2:D51").Style = "Normal"
